# Resumen de partidos hasta el dia 27 de junio
# Enter the results for the four matches played on 2018-06-27:
#   Group E: Serbia 0 - 2 Brasil
#   Group E: Suiza 1 - 1 Costa Rica
#   Group F: República de Corea 2 - 0 Alemania
#   Group F: México 0 - 3 Suecia

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2018 World Cup")

$ws.Range("F47").Value = 0
$ws.Range("G47").Value = 2

$ws.Range("F48").Value = 1
$ws.Range("G48").Value = 1

$ws.Range("F49").Value = 2
$ws.Range("G49").Value = 0

$ws.Range("F50").Value = 0
$ws.Range("G50").Value = 3

$ws.Range("G49").Select()
